$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate row 2 (A2:J2) with the new player name entries.
$ws.Range("A2").Value = "LEOOW"
$ws.Range("B2").Value = "RH77"
$ws.Range("C2").Value = "BIGGIECHEESE"
$ws.Range("D2").Value = "SEPPUNII"
$ws.Range("E2").Value = "EMPTY"
$ws.Range("F2").Value = "OPTIMALSHOT"
$ws.Range("G2").Value = "TROOG"
$ws.Range("H2").Value = "YAMISGEY"
$ws.Range("I2").Value = "ALCHEMIST"
$ws.Range("J2").Value = "EMPTY"

# P23 was stored as text "2636"; convert it to the numeric value 2636.
$ws.Range("P23").Value = 2636
